# Updates cryptos list with refreshed price/volume data (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.569.60"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.573.39"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "2.585.49"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "3.028.75"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "58.569.22"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.543.91"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "336.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -9.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.583"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "1.971.63"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
